# Edit script: applies two changes described in the diff
#  1. Remove the "_GoBack" bookmark from the "Dibagi menjadi 2 : ... manfaat bulanan" paragraph
#  2. Turn the last (trailing, empty) paragraph into a "Catatan Pengembangan :" section with
#     additional notes, variable-category legend paragraphs, a relocated "_GoBack" bookmark,
#     and a final empty paragraph before the section break.

$d = $word.ActiveDocument

# --- Part 1: remove the "_GoBack" bookmark from the middle of the document ---
$found = $d.Content.Find
$found.ClearFormatting()
$ok = $found.Execute("manfaat bulanan, dan manfaat sekaligus")
if (-not $ok) {
    throw "Could not locate target paragraph for bookmark removal"
}
$para1 = $found.Parent.Paragraphs(1)
$r1 = $para1.Range

$bookmarkFixXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008C0BF3" w:rsidRDefault="008B1DF0" w:rsidP="00847AC5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="31"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="317" w:hanging="317"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>Dibagi menjadi 2</w:t></w:r><w:r w:rsidR="00F126A0"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> :</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> manfaat bulanan, dan manfaat sekaligus.</w:t></w:r></w:p>
"@
$r1.InsertXML($bookmarkFixXml)

# --- Part 2: replace the last (trailing empty) paragraph with the new content block ---
$n = $d.Paragraphs.Count
$para2 = $d.Paragraphs($n)
$r2 = $para2.Range

$trailingXml = @"
<w:p w:rsidR="0023273F" w:rsidRPr="007D02FC" w:rsidRDefault="0023273F" w:rsidP="007D02FC"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>Catatan Pengembangan :</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>By Iqbal</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>Variabel kategori tanggungan :</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">M1 : Laki-laki menikah dan memiliki beberapa anak </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>M2 : Laki-laki menikah belum memiliki anak</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">M3 : Lajang/Duda dan memiliki beberapa anak </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>M4 : Laki-laki Lajang</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">1 : </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>Perempuan</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> menikah dan memiliki beberapa anak </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> : Perempuan</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> menikah belum memiliki anak</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">3 : </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>Lajang/Janda</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> dan memiliki beberapa anak </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">4 : </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">Perempuan </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:sz w:val="22"/><w:lang w:val="id-ID"/></w:rPr><w:t>Lajang</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p>
"@
$r2.InsertXML($trailingXml)

Write-Host "Edit complete"
